# Mise à jour des résultats du script
# Appends the latest scraped rows (2025-09-11) to the results table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append right after the current last row (238).
$newRows = @(
    @{ Row = 239; Date = "2025-09-11"; Terme = "ruissellement";          Page = 37; Occ = 2 },
    @{ Row = 240; Date = "2025-09-11"; Terme = "développement durable";  Page = 38; Occ = 1 },
    @{ Row = 241; Date = "2025-09-11"; Terme = "ruissellement";          Page = 40; Occ = 1 },
    @{ Row = 242; Date = "2025-09-11"; Terme = "ruissellement";          Page = 41; Occ = 1 }
)

# Force column A to Text so the "yyyy-mm-dd" strings are stored verbatim
# instead of being auto-converted to date serials (matches the other rows
# in the sheet, which hold the date as plain text).
$ws.Range("A239:A242").NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Date
    $ws.Range("B$row").Value = $r.Terme
    $ws.Range("C$row").Value = $r.Page
    $ws.Range("D$row").Value = $r.Occ
}
